# Apply strike-through formatting to two rubric line items.
#
# 1) "3.  absolute and relative positioning in at least one area of your
#    web page" - add strike-through to the whole line and remove its
#    yellow highlight (the line is being marked as done/obsolete).
# 2) "1.  Include JavaScript code that verifies the email fields match..."
#    - add strike-through, but keep its existing yellow highlight.

$d = $word.ActiveDocument

function Get-ParagraphByText($doc, [string]$needle) {
    $count = $doc.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        if ($p.Range.Text -like "*$needle*") {
            return $p
        }
    }
    return $null
}

# --- 1) absolute/relative positioning line: strike + drop highlight ---
$p1 = Get-ParagraphByText $d "absolute and relative positioning"
$r1 = $p1.Range
$r1.Font.StrikeThrough = 1
$r1.HighlightColorIndex = 0

# --- 2) JavaScript email validation line: strike, keep highlight ---
$p2 = Get-ParagraphByText $d "Include JavaScript code that verifies"
$r2 = $p2.Range
$r2.Font.StrikeThrough = 1
